$wb = $excel.ActiveWorkbook
$wsVerify = $wb.Worksheets.Item("VerifyOtp")
$wsGetOtp = $wb.Worksheets.Item("GetOtp")

# Clear the "null" text from A6 and B11 on the VerifyOtp sheet
$wsVerify.Range("A6").ClearContents()
$wsVerify.Range("B11").ClearContents()

# Fill in the new row of data (A12/B12) on the VerifyOtp sheet
$wsVerify.Range("A12").Value = 155
$wsVerify.Range("B12").Value = 1234

# Update the GetOtp sheet's remembered selection first (it is currently active)
[void]$wsGetOtp.Range("B11").Select()

# Make VerifyOtp the active sheet (was GetOtp) and update its selection
[void]$wsVerify.Activate()
[void]$wsVerify.Range("C16").Select()
